# Apply the diff: update OS description text, refresh pricing URLs/values
# for rows 2-4, and append a new row 5 with pricing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $addr, $val) {
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

# --- Row 2 ---
Set-TextCell $ws "B2" "Free: Debian, CentOS, CoreOS, Ubuntu or BYOL"
Set-TextCell $ws "E2" "https://cloud.google.com/products/calculator?dl=CjhDaVJpTjJSbU4ySTRaQzAyTmpJNExUUTJPREV0WW1JNU55MDNaV014TnpOak56ZzROVGNRQVE9PRAIGiQwMDRCNkZENy1GNjlFLTQxQ0UtOEI3My04NjVGNUYwQUZFMjQ"
Set-TextCell $ws "F2" "`$7.34"
Set-TextCell $ws "G2" "https://cloud.google.com/products/calculator?dl=CjhDaVJpTjJSbU4ySTRaQzAyTmpJNExUUTJPREV0WW1JNU55MDNaV014TnpOak56ZzROVGNRQVE9PRAIGiQwMDRCNkZENy1GNjlFLTQxQ0UtOEI3My04NjVGNUYwQUZFMjQ"
Set-TextCell $ws "H2" "`$7.34"
Set-TextCell $ws "I2" "https://cloud.google.com/products/calculator?dl=CjhDaVJpTjJSbU4ySTRaQzAyTmpJNExUUTJPREV0WW1JNU55MDNaV014TnpOak56ZzROVGNRQVE9PRAIGiQwMDRCNkZENy1GNjlFLTQxQ0UtOEI3My04NjVGNUYwQUZFMjQ"
Set-TextCell $ws "J2" "`$7.34"
Set-TextCell $ws "K2" "https://cloud.google.com/products/calculator?dl=CjhDaVJpTjJSbU4ySTRaQzAyTmpJNExUUTJPREV0WW1JNU55MDNaV014TnpOak56ZzROVGNRQVE9PRAIGiQwMDRCNkZENy1GNjlFLTQxQ0UtOEI3My04NjVGNUYwQUZFMjQ"
Set-TextCell $ws "L2" "`$7.34"

# --- Row 3 ---
Set-TextCell $ws "B3" "Free: Debian, CentOS, CoreOS, Ubuntu or BYOL"
Set-TextCell $ws "E3" "https://cloud.google.com/products/calculator?dl=CjhDaVExWVRSaFlqSmxOeTFqWVRaaUxUUmpaVE10WVdabU15MW1NVFl4TkROa05XSXhPRFlRQVE9PRAIGiQwRjk1RDcxNi00NTA5LTRGN0QtOTJERS1EMjY1NTc4Nzk1QUU"
Set-TextCell $ws "F3" "`$41.65"
Set-TextCell $ws "G3" "https://cloud.google.com/products/calculator?dl=CjhDaVJsWW1ObU5qQTFZeTB5TXpkbExUUmpPVEF0WVdNeU15MDFOVFk0T1RkaVpUazNaVGtRQVE9PRAIGiQwREU0QjBCMi05ODE3LTRFMDQtOUY0My0wRUNDNkZBQjg5MEU"
Set-TextCell $ws "H3" "`$29.15"
Set-TextCell $ws "I3" "https://cloud.google.com/products/calculator?dl=CjhDaVJtTnpNNU9XSTBOaTAxT1dKbUxUUXpObU10WWpVellTMWxOVEV4TkRNNFpUVXdaV1lRQVE9PRAIGiRCOUM5N0RCRS04RUM3LTQyRTctODg0Qi1DQ0E1RUFEQTY2QTU"
Set-TextCell $ws "J3" "`$26.23"
Set-TextCell $ws "K3" "https://cloud.google.com/products/calculator?dl=CjhDaVExWkdSaE1qYzBOaTB3WkRKaUxUUTFObU10T0RsaE15MDROekU0Wm1NeE0yRmpNR0lRQVE9PRAIGiRBRkVEODAyQi03MDlBLTQ3QTMtOEY0Qi0xMTZEREI5NDFBRkQ"
Set-TextCell $ws "L3" "`$18.74"

# --- Row 4 ---
Set-TextCell $ws "B4" "Free: Debian, CentOS, CoreOS, Ubuntu or BYOL"
Set-TextCell $ws "E4" "https://cloud.google.com/products/calculator?dl=CjhDaVEzT1RsalpHWmpOeTB6TnpJM0xUUTJaak10T1dKbU1pMDBOalprWW1VeE16YzJZemdRQVE9PRAIGiRFOTc0RTY5MC1COEM3LTRBMTQtODg1NS02M0YzMDVERjEwMzk"
Set-TextCell $ws "F4" "`$401.04"
Set-TextCell $ws "G4" "https://cloud.google.com/products/calculator?dl=CjhDaVF5TmpkaVptUmhOQzAwTldaakxUUTBPVFV0T0Roa01TMWhaR1ZtTnpRNE1qSTVNVElRQVE9PRAIGiRGMjhDMjVGMS1FQzk0LTQzNEItODAyNi03MjcxQzBBMTAzRTI"
Set-TextCell $ws "H4" "`$319.91"
Set-TextCell $ws "I4" "https://cloud.google.com/products/calculator?dl=CjhDaVExWTJZek1XUTJOQzFrWVRNekxUUTBNV0l0T0RRd1pTMDJaakExT1RRellqRTFNakVRQVE9PRAIGiRGQUQ5OTg4MC02ODM3LTRGRUYtOUMwMy02NTJGOTZDODlDNTA"
Set-TextCell $ws "J4" "`$254.17"
Set-TextCell $ws "K4" "https://cloud.google.com/products/calculator?dl=CjhDaVEyT1RJd01XWmtOUzFpWlRZeExUUTJNV0V0WVdRM09TMHdZekJoTkdZME1UUTJaRFlRQVE9PRAIGiQzRUIyQUUzMC0yNzg1LTQ4MEQtQjM2NC01QkY3MDMzMkZCQUE"
Set-TextCell $ws "L4" "`$185.42"

# --- Row 5 (new) ---
$ws.Range("A5").Value = 4
Set-TextCell $ws "B5" "Free: Debian, CentOS, CoreOS, Ubuntu or BYOL"
$ws.Range("C5").Value = 1
Set-TextCell $ws "D5" "general purpose"
Set-TextCell $ws "E5" "https://cloud.google.com/products/calculator?dl=CjhDaVE1WWpBNFlXSmpOUzB4WmpCakxUUmlPR0l0T1RNMllpMWlaV1JoWkRGak9EVmhOREVRQVE9PRAIGiQyM0FGODdGNi1GQUNDLTQwMDgtODAwNy04RjNDNzc3Mjc0MzA"
Set-TextCell $ws "F5" "`$109.84"
Set-TextCell $ws "G5" "https://cloud.google.com/products/calculator?dl=CjhDaVF5Wm1JNFpqUTRaUzFsTXpRd0xUUTFPR1V0T0dFME5pMDJZMkl5Wmpka1pUZ3hOVFFRQVE9PRAIGiQ0Qzk5QjExNy1CNERFLTQ2OTgtQkFCNy1ENDlDOTcyMkUzRUQ"
Set-TextCell $ws "H5" "`$87.62"
Set-TextCell $ws "I5" "Error"
Set-TextCell $ws "J5" "Error"
Set-TextCell $ws "K5" "https://cloud.google.com/products/calculator?dl=CjhDaVJqWmpoaFpEbGtNaTAzTWpBeUxUUTNPRE10WVRJNFlTMHlZamsxTTJGbE1tUTROR01RQVE9PRAIGiQzRTcyNTZGOC03MzU5LTQ3OUQtQTU4NS01REM1Qjc0MjA1OEY"
Set-TextCell $ws "L5" "`$24.71"

Write-Host "Edit complete"
